$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value2 = 0.8
$ws.Range("C2").Value2 = 0.7
$ws.Range("D2").Value2 = 0.78
$ws.Range("E2").Value2 = 0.93
$ws.Range("F2").Value2 = 0.96
$ws.Range("H2").Value2 = 0.92
$ws.Range("J2").Value2 = 0.89
$ws.Range("K2").Value2 = 0.87
# Row 3
$ws.Range("B3").Value2 = 0.74
$ws.Range("C3").Value2 = 0.67
$ws.Range("D3").Value2 = 0.76
$ws.Range("E3").Value2 = 0.97
$ws.Range("G3").Value2 = 0.89
$ws.Range("H3").Value2 = 0.93
$ws.Range("I3").Value2 = 0.93
$ws.Range("J3").Value2 = 0.88
# Row 4
$ws.Range("B4").Value2 = 0.7
$ws.Range("C4").Value2 = 0.66
$ws.Range("D4").Value2 = 0.77
$ws.Range("E4").Value2 = 0.9
$ws.Range("F4").Value2 = 0.96
$ws.Range("G4").Value2 = 0.84
$ws.Range("H4").Value2 = 0.93
$ws.Range("I4").Value2 = 0.93
$ws.Range("J4").Value2 = 0.88
$ws.Range("K4").Value2 = 0.84
# Row 5
$ws.Range("B5").Value2 = 0.77
$ws.Range("E5").Value2 = 0.97
$ws.Range("F5").Value2 = 0.99
$ws.Range("I5").Value2 = 0.95
$ws.Range("K5").Value2 = 0.87
# Row 6
$ws.Range("B6").Value2 = 0.72
$ws.Range("C6").Value2 = 0.7
$ws.Range("D6").Value2 = 0.76
$ws.Range("E6").Value2 = 0.97
$ws.Range("F6").Value2 = 0.99
$ws.Range("I6").Value2 = 0.93
$ws.Range("J6").Value2 = 0.86
$ws.Range("K6").Value2 = 0.86
# Row 7
$ws.Range("B7").Value2 = 0.78
$ws.Range("C7").Value2 = 0.64
$ws.Range("D7").Value2 = 0.78
$ws.Range("E7").Value2 = 0.92
$ws.Range("F7").Value2 = 0.98
$ws.Range("G7").Value2 = 0.91
$ws.Range("H7").Value2 = 0.92
$ws.Range("I7").Value2 = 0.9399999999999999
$ws.Range("J7").Value2 = 0.92
# Row 8
$ws.Range("C8").Value2 = 0.72
$ws.Range("D8").Value2 = 0.77
$ws.Range("E8").Value2 = 0.87
$ws.Range("G8").Value2 = 0.88
$ws.Range("J8").Value2 = 0.87
$ws.Range("K8").Value2 = 0.86
# Row 9
$ws.Range("B9").Value2 = 0.75
$ws.Range("C9").Value2 = 0.73
$ws.Range("D9").Value2 = 0.78
$ws.Range("E9").Value2 = 0.95
$ws.Range("F9").Value2 = 0.99
$ws.Range("I9").Value2 = 0.93
$ws.Range("K9").Value2 = 0.87
# Row 10
$ws.Range("C10").Value2 = 0.68
$ws.Range("D10").Value2 = 0.76
$ws.Range("E10").Value2 = 0.88
$ws.Range("G10").Value2 = 0.9
$ws.Range("H10").Value2 = 0.89
$ws.Range("J10").Value2 = 0.9
$ws.Range("K10").Value2 = 0.85
# Row 11
$ws.Range("B11").Value2 = 0.8
$ws.Range("C11").Value2 = 0.68
$ws.Range("D11").Value2 = 0.71
$ws.Range("E11").Value2 = 0.78
$ws.Range("F11").Value2 = 0.98
$ws.Range("G11").Value2 = 0.91
$ws.Range("H11").Value2 = 0.88
$ws.Range("I11").Value2 = 0.9399999999999999
$ws.Range("J11").Value2 = 0.89
$ws.Range("K11").Value2 = 0.84
# Row 12
$ws.Range("B12").Value2 = 0.79
$ws.Range("C12").Value2 = 0.65
$ws.Range("D12").Value2 = 0.74
$ws.Range("E12").Value2 = 0.85
$ws.Range("G12").Value2 = 0.91
$ws.Range("H12").Value2 = 0.88
$ws.Range("I12").Value2 = 0.96
$ws.Range("J12").Value2 = 0.88
# Row 13
$ws.Range("B13").Value2 = 0.78
$ws.Range("C13").Value2 = 0.71
$ws.Range("D13").Value2 = 0.77
$ws.Range("E13").Value2 = 0.9399999999999999
$ws.Range("F13").Value2 = 0.97
$ws.Range("G13").Value2 = 0.89
$ws.Range("H13").Value2 = 0.91
$ws.Range("I13").Value2 = 0.9399999999999999
$ws.Range("J13").Value2 = 0.91
$ws.Range("K13").Value2 = 0.87
# Row 14
$ws.Range("B14").Value2 = 0.8
$ws.Range("C14").Value2 = 0.68
$ws.Range("E14").Value2 = 0.96
$ws.Range("F14").Value2 = 0.97
$ws.Range("G14").Value2 = 0.8100000000000001
$ws.Range("H14").Value2 = 0.9
$ws.Range("J14").Value2 = 0.88
$ws.Range("K14").Value2 = 0.86
# Row 15
$ws.Range("B15").Value2 = 0.78
$ws.Range("C15").Value2 = 0.7
$ws.Range("E15").Value2 = 0.9399999999999999
$ws.Range("F15").Value2 = 0.98
$ws.Range("G15").Value2 = 0.88
$ws.Range("K15").Value2 = 0.86
# Row 16
$ws.Range("F16").Value2 = 0.98
